# ---------------------------------------------------------------------------
# POM registration/assert with all my profile fields
#
# This script:
#   1. Inserts a new "country" column into the registrationInfo sheet and
#      refreshes the registration / phone / e-mail data for the three users.
#   2. Adds a new "verifyRegistration" worksheet that captures the
#      email/password/country used to verify each registration.
#   3. Adds a blank "Sheet2" placeholder worksheet at the end of the workbook.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. registrationInfo sheet
# ---------------------------------------------------------------------------
$reg = $wb.Worksheets.Item("registrationInfo")
$reg.Activate()

# Insert a new column before the old "Address" column (F) for "country".
$reg.Columns.Item(6).Insert()

# Header
$reg.Range("F1").Value = "country"

# Drop the existing (stale) hyperlinks on the Email column before we
# overwrite the addresses - they will be re-created below.
$reg.Hyperlinks.Delete()

# Row 2 - Masha L Rey
$reg.Range("D2").Value = "masharey915677@yahoo.com"
$reg.Range("E2").Value = "(789) 876-7877"
$reg.Range("F2").Value = "Latvia"

# Row 3 - John M Smith
$reg.Range("D3").Value = "johnjmith345677@gmail.com"
$reg.Range("E3").Value = "(789) 898-7890"
$reg.Range("F3").Value = "Serbia"
$reg.Range("H3").Value = "Apt.2"

# Row 4 - Alex L Erny
$reg.Range("D4").Value = "alexerny345677@gmail.com"
$reg.Range("E4").Value = "(676) 567-8767"
$reg.Range("F4").Value = "Romania"

# Re-create the hyperlinks for the e-mail column, then re-apply the
# "Hyperlink" cell style (Add() on its own drops the original styling).
$reg.Hyperlinks.Add($reg.Range("D2"), "mailto:masharey915677@yahoo.com") | Out-Null
$reg.Hyperlinks.Add($reg.Range("D3"), "mailto:johnjmith345677@gmail.com") | Out-Null
$reg.Hyperlinks.Add($reg.Range("D4"), "mailto:alexerny345677@gmail.com") | Out-Null
$reg.Range("D2").Style = "Hyperlink"
$reg.Range("D3").Style = "Hyperlink"
$reg.Range("D4").Style = "Hyperlink"

# Give the new "country" column a fixed (non best-fit) width of 11.
$reg.Columns.Item(6).ColumnWidth = 10.1

# Re-fit the phone column now that the values are formatted strings.
$reg.Columns.Item(5).EntireColumn.AutoFit() | Out-Null

$reg.Range("D3").Select() | Out-Null

# ---------------------------------------------------------------------------
# 2. verifyRegistration sheet (new)
# ---------------------------------------------------------------------------
$afterSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$verify = $wb.Worksheets.Add($null, $afterSheet)
$verify.Name = "verifyRegistration"

$verify.Range("A1").Value = "email"
$verify.Range("B1").Value = "pass"
$verify.Range("C1").Value = "country"

$verify.Range("A2").Value = "AlexErny21@gmail.com"
$verify.Range("B2").Value = "tampaflorida"
$verify.Range("C2").Value = "Romania"

$verify.Range("A3").Value = "JohnSmith12@gmail.com"
$verify.Range("B3").Value = "John12345"
$verify.Range("C3").Value = "Serbia"

$verify.Range("A4").Value = "MashaRey12345@yahoo.com"
$verify.Range("B4").Value = "masha234"
$verify.Range("C4").Value = "Latvia"

$verify.Hyperlinks.Add($verify.Range("A4"), "mailto:MashaRey12345@yahoo.com") | Out-Null
$verify.Range("A4").Style = "Hyperlink"

$verify.Columns.Item(1).EntireColumn.AutoFit() | Out-Null
$verify.Columns.Item(2).EntireColumn.AutoFit() | Out-Null

$verify.Range("A4").Select() | Out-Null

# ---------------------------------------------------------------------------
# 3. Sheet2 (new, blank placeholder)
# ---------------------------------------------------------------------------
$afterSheet2 = $wb.Worksheets.Item($wb.Worksheets.Count)
$blank = $wb.Worksheets.Add($null, $afterSheet2)
$blank.Name = "Sheet2"

# ---------------------------------------------------------------------------
# Re-activate the registration sheet so it stays the selected tab.
# ---------------------------------------------------------------------------
$reg.Activate()
